$wb = $excel.ActiveWorkbook

# --- Fix typo "peroper" -> "proper" in the "links" sheet (column A, rows 2-17) ---
$links = $wb.Worksheets.Item("links")
$lastLinkRow = 17
for ($r = 2; $r -le $lastLinkRow; $r++) {
    $cell = $links.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null -and $old.Contains("peroper")) {
        $cell.Value = $old.Replace("peroper", "proper")
    }
}

# --- Fix wording "asserting article name" -> "assert article name" in "relatedArticles" sheet (column A, rows 2-12) ---
$relatedArticles = $wb.Worksheets.Item("relatedArticles")
$lastArticleRow = 12
for ($r = 2; $r -le $lastArticleRow; $r++) {
    $cell = $relatedArticles.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null -and $old.Contains("asserting article name")) {
        $cell.Value = $old.Replace("asserting article name", "assert article name")
    }
}

# --- Update selections left behind on "links" sheet ---
[void]$links.Activate()
[void]$links.Range("A17").Select()

# --- Make "relatedArticles" the active sheet/tab, with its own lingering selection ---
[void]$relatedArticles.Activate()
[void]$relatedArticles.Range("A12").Select()
